# Apply the updated crypto price/volume snapshot (GitHub Actions data refresh).
# Source data cells are plain text (inlineStr) in the workbook: numeric-looking
# price strings (column D) must be forced to Text format before assignment so the
# COM layer does not silently coerce them into numbers (which would drop meaningful
# trailing/leading zeros, e.g. "152.10" -> 152.1 or "0.110" -> 0.11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '62.206.68'
$ws.Cells.Item(2, 5).Value = '  +0.70%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '2.419.60'
$ws.Cells.Item(3, 5).Value = '  +0.90%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5: BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(5, 4).Value = '564.15'
$ws.Cells.Item(5, 5).Value = '  +2.08%  '

# Row 6: Solana
$ws.Cells.Item(6, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(6, 4).Value = '143.05'
$ws.Cells.Item(6, 5).Value = '  +0.59%  '

# Row 8: XRP
$ws.Cells.Item(8, 5).Value = '  +1.91%  '

# Row 9: LidoStakedEther
$ws.Cells.Item(9, 4).Value = '2.415.95'
$ws.Cells.Item(9, 5).Value = '  +0.95%  '

# Row 10: Dogecoin
$ws.Cells.Item(10, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(10, 4).Value = '0.110'
$ws.Cells.Item(10, 5).Value = '  +1.81%  '

# Row 11: TRON
$ws.Cells.Item(11, 5).Value = '  -2.10%  '

# Row 12: Toncoin
$ws.Cells.Item(12, 5).Value = '  -0.40%  '

# Row 13: Cardano
$ws.Cells.Item(13, 5).Value = '  +0.36%  '

# Row 14: Avalanche
$ws.Cells.Item(14, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(14, 4).Value = '25.75'
$ws.Cells.Item(14, 5).Value = '  -1.00%  '

# Row 15: ShibaInu
$ws.Cells.Item(15, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(15, 4).Value = '0.0000176'
$ws.Cells.Item(15, 5).Value = '  +0.82%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Cells.Item(16, 4).Value = '2.854.93'
$ws.Cells.Item(16, 5).Value = '  +1.03%  '

# Row 17: WrappedBTC
$ws.Cells.Item(17, 4).Value = '61.985.81'
$ws.Cells.Item(17, 5).Value = '  +0.62%  '

# Row 18: WrappedEther
$ws.Cells.Item(18, 4).Value = '2.415.69'
$ws.Cells.Item(18, 5).Value = '  +1.05%  '

# Row 19: Chainlink
$ws.Cells.Item(19, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(19, 4).Value = '11.35'
$ws.Cells.Item(19, 5).Value = '  +2.00%  '

# Row 20: Polkadot -> BitcoinCash
$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(20, 4).Value = '323.83'
$ws.Cells.Item(20, 5).Value = '  +0.51%  '

# Row 21: BitcoinCash -> Polkadot
$ws.Cells.Item(21, 2).Value = 'Polkadot'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(21, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(21, 4).Value = '4.18'
$ws.Cells.Item(21, 5).Value = '  +0.42%  '

# Row 22: Uniswap
$ws.Cells.Item(22, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(22, 4).Value = '6.86'
$ws.Cells.Item(22, 5).Value = '  +2.94%  '

# Row 23: Dai
$ws.Cells.Item(23, 5).Value = '  -0.15%  '

# Row 24: Litecoin
$ws.Cells.Item(24, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(24, 4).Value = '65.98'
$ws.Cells.Item(24, 5).Value = '  +2.32%  '

# Row 26: Aptos
$ws.Cells.Item(26, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(26, 4).Value = '8.96'
$ws.Cells.Item(26, 5).Value = '  -3.01%  '

# Row 27: Bittensor
$ws.Cells.Item(27, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(27, 4).Value = '581.28'
$ws.Cells.Item(27, 5).Value = '  +4.00%  '

# Row 28: PEPE
$ws.Cells.Item(28, 4).Value = '0.0₃0950'
$ws.Cells.Item(28, 5).Value = '  +3.45%  '

# Row 29: WrappedeETH
$ws.Cells.Item(29, 4).Value = '2.531.19'

# Row 30: Binance-PegBSC-USD
$ws.Cells.Item(30, 5).Value = '  -0.56%  '

# Row 31: InternetComputer(DFINITY)
$ws.Cells.Item(31, 5).Value = '  -0.46%  '

# Row 32: Fetch.AI
$ws.Cells.Item(32, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(32, 4).Value = '1.44'
$ws.Cells.Item(32, 5).Value = '  +1.04%  '

# Row 33: Kaspa
$ws.Cells.Item(33, 5).Value = '  +0.40%  '

# Row 34: PancakeSwap
$ws.Cells.Item(34, 5).Value = '  +0.99%  '

# Row 35: ImmutableX
$ws.Cells.Item(35, 5).Value = '  +0.20%  '

# Row 36: FirstDigitalUSD
$ws.Cells.Item(36, 5).Value = '  +0.28%  '

# Row 37: RenderToken
$ws.Cells.Item(37, 5).Value = '  -3.67%  '

# Row 38: NEARProtocol
$ws.Cells.Item(38, 5).Value = '  -0.20%  '

# Row 40: Monero
$ws.Cells.Item(40, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(40, 4).Value = '152.03'
$ws.Cells.Item(40, 5).Value = '  +4.12%  '

# Row 41: EthereumClassic
$ws.Cells.Item(41, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(41, 4).Value = '18.70'
$ws.Cells.Item(41, 5).Value = '  +0.42%  '

# Row 42: Stacks
$ws.Cells.Item(42, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(42, 4).Value = '1.80'
$ws.Cells.Item(42, 5).Value = '  -7.86%  '

# Row 43: USDe
$ws.Cells.Item(43, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  -0.07%  '

# Row 44: dogwifhat
$ws.Cells.Item(44, 5).Value = '  +1.10%  '

# Row 45: Aave
$ws.Cells.Item(45, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(45, 4).Value = '148.66'
$ws.Cells.Item(45, 5).Value = '  +0.01%  '

# Row 46: Filecoin
$ws.Cells.Item(46, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(46, 4).Value = '3.66'
$ws.Cells.Item(46, 5).Value = '  +0.67%  '

# Row 47: Hedera
$ws.Cells.Item(47, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(47, 4).Value = '0.0536'
$ws.Cells.Item(47, 5).Value = '  +1.09%  '

# Row 48: InjectiveProtocol
$ws.Cells.Item(48, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(48, 4).Value = '20.10'
$ws.Cells.Item(48, 5).Value = '  -0.72%  '

# Row 49: Mantle
$ws.Cells.Item(49, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(49, 4).Value = '0.595'
$ws.Cells.Item(49, 5).Value = '  +1.60%  '

# Row 50: Stellar
$ws.Cells.Item(50, 4).NumberFormat = "@"  # keep price as text
$ws.Cells.Item(50, 4).Value = '0.0920'
$ws.Cells.Item(50, 5).Value = '  +1.67%  '

# Row 51: VeChain
$ws.Cells.Item(51, 5).Value = '  +1.49%  '
